$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$pic = $ws.Shapes.AddPicture("/tmp/work/extracted/xl/media/image1.png", $false, $true, 700, 300, 20, 23)
Write-Host "Added: $($pic.Name)"
Write-Host "Count: $($ws.Shapes.Count)"
